# Coastal Surface Piercing Profilers - update Omaha Cal Info for
# CP05MOAS-GL001 -> CP05MOAS-GL374 (Mooring + Asset_Cal_Info sheets)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # Moorings
$ws2 = $wb.Worksheets.Item(2)  # Asset_Cal_Info

# --- Update the Mooring Serial Number on the "Moorings" sheet ---
$ws1.Range("A2").Value = "CP05MOAS-GL374"

# --- Update the Ref Des values on the "Asset_Cal_Info" sheet ---
# ADCPAM000 block
$ws2.Range("A2").Value = "CP05MOAS-GL374-01-ADCPAM000"
$ws2.Range("A3").Value = "CP05MOAS-GL374-01-ADCPAM000"
$ws2.Range("A4").Value = "CP05MOAS-GL374-01-ADCPAM000"
$ws2.Range("A5").Value = "CP05MOAS-GL374-01-ADCPAM000"

# FLORTM000 block
$ws2.Range("A7").Value = "CP05MOAS-GL374-02-FLORTM000"
$ws2.Range("A8").Value = "CP05MOAS-GL374-02-FLORTM000"
$ws2.Range("A9").Value = "CP05MOAS-GL374-02-FLORTM000"
$ws2.Range("A10").Value = "CP05MOAS-GL374-02-FLORTM000"

# CTDGVM000
$ws2.Range("A12").Value = "CP05MOAS-GL374-03-CTDGVM000"

# DOSTAM000
$ws2.Range("A14").Value = "CP05MOAS-GL374-04-DOSTAM000"

# PARADM000
$ws2.Range("A16").Value = "CP05MOAS-GL374-05-PARADM000"

# ENG000000
$ws2.Range("A18").Value = "CP05MOAS-GL374-00-ENG000000"

# --- Update the active sheet / selections to match the saved UI state ---
# Moorings: no longer the selected tab, but its remembered selection moves to D8
[void]$ws1.Activate()
[void]$ws1.Range("D8").Select()

# Asset_Cal_Info becomes the selected/active tab, selection stays at F9
[void]$ws2.Activate()
[void]$ws2.Range("F9").Select()
